# "list to array issue" - convert the Plagiarised column (B) from sparse
# True/False/blank booleans into a fully populated 0/1 numeric array.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New numeric values (0/1) for B2:B55, replacing the old boolean/blank cells.
$values = @(
    1,1,1,1,1,1,1,1,1,1,1,1,   # rows 2-13
    0,0,0,0,0,0,0,             # rows 14-20
    1,                          # row 21
    0,0,0,0,0,0,0,0,0,0,0,     # rows 22-32
    1,1,1,1,1,1,1,1,1,1,1,1,1,1, # rows 33-46
    0,0,0,0,0,0,0,0,0          # rows 47-55
)

$startRow = 2
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update the view state left in the sheet after the edit: scrolled down
# so row 16 is at the top, with B47:B55 selected (active cell B47).
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("B47:B55").Select()
